$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 8 with the new "Isotonic" calibration method results
$ws.Range("A8").Value = "Classification"
$ws.Range("B8").Value = "CIFAR100"
$ws.Range("C8").Value = "DenseNet"
$ws.Range("D8").Value = "Isotonic"
$ws.Range("F8").Value = 0.74968000000000001
$ws.Range("G8").Value = 0.065079999999999999
$ws.Range("H8").Value = 0.047219999999999998
$ws.Range("I8").Value = 0.1356

# Copy styles from row 7 (A7:H7) down to row 8 (A8:H8), matching the bordered style
$ws.Range("A7:H7").Copy()
$ws.Range("A8:H8").PasteSpecial(-4122)  # xlPasteFormats

# I8 should use the same style as I7 (border + number format)
$ws.Range("I7").Copy()
$ws.Range("I8").PasteSpecial(-4122)  # xlPasteFormats

# Update the active selection to D9, matching the diff
$ws.Range("D9").Select()
